$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as literal text in the source data
# (e.g. "63.549.85", "0.999", "167.70") -- force text via NumberFormat "@"
# before assigning so Excel does not reinterpret number-looking strings as
# numeric values (which would also lose trailing zeros / introduce float error).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.629.21'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.613.69'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '595.95'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.43'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.590'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.69'
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.383'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('E12').Value = '  -1.15%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.71'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.083.08'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.427.86'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.605.60'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.32'
$ws.Range('E18').Value = '  +6.22%  '
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '346.47'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.74'
$ws.Range('E23').Value = '  +3.11%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.69'
$ws.Range('E24').Value = '  -0.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.71'
$ws.Range('E25').Value = '  +8.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.21'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '552.43'
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.12'
$ws.Range('E29').Value = '  +1.64%  '
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0₃0845'
$ws.Range('E33').Value = '  -1.69%  '
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '167.69'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.413'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '166.62'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '39.72'
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('E44').Value = '  +2.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0586'
$ws.Range('E45').Value = '  +2.06%  '
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0₆0251'
$ws.Range('E49').Value = '  +24.55%  '
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0965'
$ws.Range('E51').Value = '  +0.08%  '
